# Add team record columns (Wins / Losses / Ties) to the right of the
# existing data, mirroring the "Unnamed: 28" column's header styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the formatting of the last existing header
# cell (AC1) onto the three new header cells so they pick up the same
# bold/centered/bordered style, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-54): every team/player row gets the same W/L/T record.
$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 65   # AD = Wins
    $ws.Cells.Item($r, 31).Value = 97   # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = Ties
}
